# Add a new "2021" column (R) to the worksheet, mirroring the existing
# 2020 column (Q): same formatting, with updated figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: year header -------------------------------------------------
$ws.Range("Q3").Copy()
$ws.Range("R3").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("R3").Value = 2021

# --- Row 4: a) branches per 100 000 adults (formula) ---------------------
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)
$ws.Range("R4").Formula = "=R6/R8*100000"

# --- Row 5: b) ATMs per 100 000 adults (formula) --------------------------
$ws.Range("Q5").Copy()
$ws.Range("R5").PasteSpecial(-4122)
$ws.Range("R5").Formula = "=R7/R8*100000"

# --- Row 6: total commercial bank branches --------------------------------
$ws.Range("Q6").Copy()
$ws.Range("R6").PasteSpecial(-4122)
$ws.Range("R6").Value = 312

# --- Row 7: total ATMs ----------------------------------------------------
$ws.Range("Q7").Copy()
$ws.Range("R7").PasteSpecial(-4122)
$ws.Range("R7").Value = 1910

# --- Row 8: adult resident population -------------------------------------
$ws.Range("Q8").Copy()
$ws.Range("R8").PasteSpecial(-4122)
$ws.Range("R8").Value = 4409166

$excel.CutCopyMode = $false

# Update the sheet view: move the selection (this also clears the old
# scrolled-in topLeftCell position).
$ws.Range("R15").Select()
